{"js": "// Applies the replacement of each paragraph's text (date line + 100 table-cell\n// math expressions) with its new value, in document order. The new values were\n// derived from the target OOXML diff; order corresponds 1:1 with\n// context.document.body.paragraphs (which enumerates paragraphs -- including\n// those that live inside table cells -- in document order).\nconst newValues = [\n  \"2023-08-19 Saturday\",\n  \"48+18=66\",\n  \"56-48=8\",\n  \"33+38=71\",\n  \"9+19=28\",\n  \"49+8=57\",\n  \"60-46=14\",\n  \"3+89=92\",\n  \"54-49=5\",\n  \"58+5=63\",\n  \"43-36=7\",\n  \"65+26=91\",\n  \"7+64=71\",\n  \"57+28=85\",\n  \"32-19=13\",\n  \"65-9=56\",\n  \"74-55=19\",\n  \"8+88=96\",\n  \"71-32=39\",\n  \"39+14=53\",\n  \"94-6=88\",\n  \"6+65=71\",\n  \"91-48=43\",\n  \"80-53=27\",\n  \"67+25=92\",\n  \"14+9=23\",\n  \"35+48=83\",\n  \"76+7=83\",\n  \"38+56=94\",\n  \"9+18=27\",\n  \"9+36=45\",\n  \"93-17=76\",\n  \"20-17=3\",\n  \"45+18=63\",\n  \"54+19=73\",\n  \"43-6=37\",\n  \"81-52=29\",\n  \"59+26=85\",\n  \"87-68=19\",\n  \"90-7=83\",\n  \"83-9=74\",\n  \"91-72=19\",\n  \"60-59=1\",\n  \"83-6=77\",\n  \"9+14=23\",\n  \"95-17=78\",\n  \"73-18=55\",\n  \"14+77=91\",\n  \"66-7=59\",\n  \"93-36=57\",\n  \"17+44=61\",\n  \"27+58=85\",\n  \"55+19=74\",\n  \"36+58=94\",\n  \"56+15=71\",\n  \"85-6=79\",\n  \"4+57=61\",\n  \"90-74=16\",\n  \"69+27=96\",\n  \"19+69=88\",\n  \"44+48=92\",\n  \"38+26=64\",\n  \"57+17=74\",\n  \"5+8=13\",\n  \"65+19=84\",\n  \"43-24=19\",\n  \"70-18=52\",\n  \"43-7=36\",\n  \"44+47=91\",\n  \"26+6=32\",\n  \"70-32=38\",\n  \"24-5=19\",\n  \"26-19=7\",\n  \"39+2=41\",\n  \"39+6=45\",\n  \"75-69=6\",\n  \"88-59=29\",\n  \"70-45=25\",\n  \"55-47=8\",\n  \"49+44=93\",\n  \"81-24=57\",\n  \"87+9=96\",\n  \"59+32=91\",\n  \"27+29=56\",\n  \"46+49=95\",\n  \"51-24=27\",\n  \"80-25=55\",\n  \"51-3=48\",\n  \"4+17=21\",\n  \"94-85=9\",\n  \"2+29=31\",\n  \"9+15=24\",\n  \"85-68=17\",\n  \"22-3=19\",\n  \"17+26=43\",\n  \"48+14=62\",\n  \"17+17=34\",\n  \"46+6=52\",\n  \"7+56=63\",\n  \"59+6=65\",\n  \"65-18=47\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newValues.length +\n    \" but found \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < newValues.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replaces the text of the date paragraph and each of the 100 table-cell\n# paragraphs (5 columns x 20 rows) with the values from the target edit, in\n# document order. Word COM's $d.Paragraphs collection also yields one empty\n# paragraph per table row (the row-end mark) plus a trailing empty paragraph\n# after the table, so we skip any paragraph whose visible text (after\n# stripping the paragraph mark / cell mark control characters) is empty.\n$newValues = @(\n    \"2023-08-19 Saturday\",\n    \"48+18=66\",\n    \"56-48=8\",\n    \"33+38=71\",\n    \"9+19=28\",\n    \"49+8=57\",\n    \"60-46=14\",\n    \"3+89=92\",\n    \"54-49=5\",\n    \"58+5=63\",\n    \"43-36=7\",\n    \"65+26=91\",\n    \"7+64=71\",\n    \"57+28=85\",\n    \"32-19=13\",\n    \"65-9=56\",\n    \"74-55=19\",\n    \"8+88=96\",\n    \"71-32=39\",\n    \"39+14=53\",\n    \"94-6=88\",\n    \"6+65=71\",\n    \"91-48=43\",\n    \"80-53=27\",\n    \"67+25=92\",\n    \"14+9=23\",\n    \"35+48=83\",\n    \"76+7=83\",\n    \"38+56=94\",\n    \"9+18=27\",\n    \"9+36=45\",\n    \"93-17=76\",\n    \"20-17=3\",\n    \"45+18=63\",\n    \"54+19=73\",\n    \"43-6=37\",\n    \"81-52=29\",\n    \"59+26=85\",\n    \"87-68=19\",\n    \"90-7=83\",\n    \"83-9=74\",\n    \"91-72=19\",\n    \"60-59=1\",\n    \"83-6=77\",\n    \"9+14=23\",\n    \"95-17=78\",\n    \"73-18=55\",\n    \"14+77=91\",\n    \"66-7=59\",\n    \"93-36=57\",\n    \"17+44=61\",\n    \"27+58=85\",\n    \"55+19=74\",\n    \"36+58=94\",\n    \"56+15=71\",\n    \"85-6=79\",\n    \"4+57=61\",\n    \"90-74=16\",\n    \"69+27=96\",\n    \"19+69=88\",\n    \"44+48=92\",\n    \"38+26=64\",\n    \"57+17=74\",\n    \"5+8=13\",\n    \"65+19=84\",\n    \"43-24=19\",\n    \"70-18=52\",\n    \"43-7=36\",\n    \"44+47=91\",\n    \"26+6=32\",\n    \"70-32=38\",\n    \"24-5=19\",\n    \"26-19=7\",\n    \"39+2=41\",\n    \"39+6=45\",\n    \"75-69=6\",\n    \"88-59=29\",\n    \"70-45=25\",\n    \"55-47=8\",\n    \"49+44=93\",\n    \"81-24=57\",\n    \"87+9=96\",\n    \"59+32=91\",\n    \"27+29=56\",\n    \"46+49=95\",\n    \"51-24=27\",\n    \"80-25=55\",\n    \"51-3=48\",\n    \"4+17=21\",\n    \"94-85=9\",\n    \"2+29=31\",\n    \"9+15=24\",\n    \"85-68=17\",\n    \"22-3=19\",\n    \"17+26=43\",\n    \"48+14=62\",\n    \"17+17=34\",\n    \"46+6=52\",\n    \"7+56=63\",\n    \"59+6=65\",\n    \"65-18=47\"\n)\n\n$d = $word.ActiveDocument\n$total = $d.Paragraphs.Count\n\n$idx = 0\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $clean = $r.Text -replace \"[`r`a]\", \"\"\n    if ($clean -ne \"\") {\n        if ($idx -ge $newValues.Length) {\n            throw \"More non-empty paragraphs than expected replacement values\"\n        }\n        $r.Text = $newValues[$idx]\n        $idx = $idx + 1\n    }\n}\n\nif ($idx -ne $newValues.Length) {\n    throw \"Expected to replace $($newValues.Length) paragraphs but replaced $idx\"\n}\n"}
